# Update cryptos list: refresh Price (col D) / Volume(1h) (col E) values for rows 2-51.
# Numeric-looking prices are prefixed with a leading apostrophe so Excel keeps
# them as text (matching the workbook's original inlineStr/text cell type)
# instead of re-interpreting them as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.125.39"
$ws.Range("E2").Value = "  -5.17%  "
$ws.Range("D3").Value = "2.994.89"
$ws.Range("E3").Value = "  -5.59%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'569.09"
$ws.Range("E5").Value = "  -4.60%  "
$ws.Range("D6").Value = "'124.88"
$ws.Range("E6").Value = "  -7.97%  "
$ws.Range("E7").Value = "  +0.23%  "
$ws.Range("D8").Value = "2.989.57"
$ws.Range("E8").Value = "  -5.75%  "
$ws.Range("E9").Value = "  -2.56%  "
$ws.Range("E10").Value = "  -7.65%  "
$ws.Range("E11").Value = "  -5.06%  "
$ws.Range("E12").Value = "  -3.08%  "
$ws.Range("E13").Value = "  -7.88%  "
$ws.Range("D14").Value = "'32.45"
$ws.Range("E14").Value = "  -6.10%  "
$ws.Range("E15").Value = "  +0.22%  "
$ws.Range("D16").Value = "3.493.96"
$ws.Range("E16").Value = "  -5.42%  "
$ws.Range("D17").Value = "2.998.69"
$ws.Range("E17").Value = "  -5.41%  "
$ws.Range("D18").Value = "60.130.35"
$ws.Range("E18").Value = "  -5.17%  "
$ws.Range("D19").Value = "'6.50"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("D20").Value = "'429.03"
$ws.Range("E20").Value = "  -7.05%  "
$ws.Range("D21").Value = "'13.13"
$ws.Range("E21").Value = "  -5.91%  "
$ws.Range("D22").Value = "'0.671"
$ws.Range("E22").Value = "  -3.72%  "
$ws.Range("D23").Value = "'7.06"
$ws.Range("E23").Value = "  -7.82%  "
$ws.Range("D24").Value = "'12.85"
$ws.Range("E24").Value = "  -2.68%  "
$ws.Range("D25").Value = "'79.43"
$ws.Range("E25").Value = "  -4.49%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'2.52"
$ws.Range("E28").Value = "  -6.15%  "
$ws.Range("E29").Value = "  -5.01%  "
$ws.Range("D30").Value = "'7.19"
$ws.Range("E30").Value = "  -7.21%  "
$ws.Range("D31").Value = "'6.11"
$ws.Range("E31").Value = "  -10.44%  "
$ws.Range("D32").Value = "'25.22"
$ws.Range("E32").Value = "  -7.44%  "
$ws.Range("D33").Value = "'0.0949"
$ws.Range("E33").Value = "  -5.93%  "
$ws.Range("D34").Value = "'5.60"
$ws.Range("E34").Value = "  -4.88%  "
$ws.Range("D35").Value = "'0.930"
$ws.Range("E35").Value = "  -9.02%  "
$ws.Range("D36").Value = "'50.32"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").Value = "'2.02"
$ws.Range("E37").Value = "  -16.26%  "
$ws.Range("D38").Value = "'8.49"
$ws.Range("E38").Value = "  +4.43%  "
$ws.Range("D39").Value = "0.0₃0660"
$ws.Range("E39").Value = "  -10.74%  "
$ws.Range("D40").Value = "'0.0356"
$ws.Range("E40").Value = "  -8.66%  "
$ws.Range("E41").Value = "  -4.22%  "
$ws.Range("D42").Value = "'370.47"
$ws.Range("E42").Value = "  -5.54%  "
$ws.Range("D43").Value = "2.673.56"
$ws.Range("E43").Value = "  -4.14%  "
$ws.Range("D44").Value = "'2.43"
$ws.Range("E44").Value = "  -7.75%  "
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "'121.26"
$ws.Range("E46").Value = "  -4.86%  "
$ws.Range("D47").Value = "'0.233"
$ws.Range("E47").Value = "  -7.05%  "
$ws.Range("D48").Value = "'1.98"
$ws.Range("E48").Value = "  -6.27%  "
$ws.Range("D49").Value = "'0.107"
$ws.Range("E49").Value = "  -3.75%  "
$ws.Range("D50").Value = "'23.24"
$ws.Range("E50").Value = "  -7.22%  "
$ws.Range("D51").Value = "'0.133"
$ws.Range("E51").Value = "  -2.20%  "
